$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 23 updates
$ws.Range("G23").Value = 1.021
$ws.Range("H23").Value = 26

# Row 29 and 40 count correction
$ws.Range("H29").Value = -1
$ws.Range("H40").Value = -1

# RF column recalculated for rows 24-41
$newRF = 41.09566265060241
for ($r = 24; $r -le 41; $r++) {
    $ws.Range("I$r").Value = $newRF
}
